# Auto-generated edit script applying numeric updates to distractor analysis tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("correct")
$ws.Range("C2").Value = 808
$ws.Range("D2").Value = 0.319
$ws.Range("C3").Value = 763
$ws.Range("D3").Value = 0.353
$ws.Range("C4").Value = 740
$ws.Range("D4").Value = 0.376
$ws.Range("C5").Value = 705
$ws.Range("D5").Value = 0.391
$ws.Range("C6").Value = 622
$ws.Range("D6").Value = 0.368
$ws.Range("C7").Value = 588
$ws.Range("D7").Value = 0.405
$ws.Range("C8").Value = 541
$ws.Range("D8").Value = 0.392
$ws.Range("C9").Value = 490
$ws.Range("D9").Value = 0.385
$ws.Range("C10").Value = 431
$ws.Range("D10").Value = 0.358
$ws.Range("C11").Value = 365
$ws.Range("D11").Value = 0.38
$ws.Range("C12").Value = 297
$ws.Range("D12").Value = 0.415
$ws.Range("C13").Value = 250
$ws.Range("D13").Value = 0.351
$ws.Range("C14").Value = 163
$ws.Range("C15").Value = 114
$ws.Range("D15").Value = 0.395
$ws.Range("C16").Value = 56
$ws.Range("D16").Value = 0.307
$ws = $wb.Worksheets.Item("distractor")
$ws.Range("C2").Value = 62
$ws.Range("D2").Value = -0.232
$ws.Range("D3").Value = -0.103
$ws.Range("C4").Value = 45
$ws.Range("D4").Value = -0.18
$ws.Range("C5").Value = 63
$ws.Range("D5").Value = -0.171
$ws.Range("C6").Value = 79
$ws.Range("D6").Value = -0.202
$ws.Range("C7").Value = 54
$ws.Range("C8").Value = 68
$ws.Range("D8").Value = -0.175
$ws.Range("C9").Value = 84
$ws.Range("D9").Value = -0.237
$ws.Range("D10").Value = -0.18
$ws.Range("C11").Value = 83
$ws.Range("D11").Value = -0.185
$ws.Range("C12").Value = 74
$ws.Range("D12").Value = -0.221
$ws.Range("C13").Value = 97
$ws.Range("D13").Value = -0.204
$ws.Range("C14").Value = 107
$ws.Range("D14").Value = -0.206
$ws.Range("C15").Value = 123
$ws.Range("D15").Value = -0.204
$ws.Range("C16").Value = 101
$ws.Range("D16").Value = -0.135
$ws.Range("C17").Value = 109
$ws.Range("D17").Value = -0.212
$ws.Range("C18").Value = 133
$ws.Range("D18").Value = -0.195
$ws.Range("C19").Value = 116
$ws.Range("D19").Value = -0.185
$ws.Range("C20").Value = 126
$ws.Range("D20").Value = -0.163
$ws.Range("C21").Value = 138
$ws.Range("D21").Value = -0.187
$ws.Range("C22").Value = 136
$ws.Range("D22").Value = -0.206
$ws.Range("C23").Value = 146
$ws.Range("D23").Value = -0.156
$ws.Range("C24").Value = 139
$ws.Range("D24").Value = -0.2
$ws.Range("C25").Value = 151
$ws.Range("D25").Value = -0.174
$ws.Range("C26").Value = 145
$ws.Range("D26").Value = -0.194
$ws.Range("C27").Value = 161
$ws.Range("D27").Value = -0.137
$ws.Range("C28").Value = 187
$ws.Range("D28").Value = -0.139
$ws.Range("C29").Value = 152
$ws.Range("D29").Value = -0.097
$ws.Range("C30").Value = 203
$ws.Range("D30").Value = -0.176
$ws.Range("C31").Value = 175
$ws.Range("D31").Value = -0.194
$ws.Range("C32").Value = 183
$ws.Range("D32").Value = -0.194
$ws.Range("C33").Value = 189
$ws.Range("D33").Value = -0.158
$ws.Range("C34").Value = 173
$ws.Range("D34").Value = -0.13
$ws.Range("C35").Value = 185
$ws.Range("D35").Value = -0.188
$ws.Range("D36").Value = -0.079
$ws.Range("C37").Value = 171
$ws.Range("D37").Value = -0.122
$ws.Range("C38").Value = 180
$ws.Range("D38").Value = -0.116
$ws.Range("C39").Value = 173
$ws.Range("D39").Value = -0.123
$ws.Range("C40").Value = 163
$ws.Range("D40").Value = -0.098
$ws.Range("C41").Value = 151
$ws.Range("D41").Value = -0.127
$ws.Range("D42").Value = -0.12
$ws.Range("C43").Value = 144
$ws.Range("D43").Value = -0.116
$ws.Range("C44").Value = 98
$ws.Range("D44").Value = -0.04
$ws.Range("C45").Value = 97
$ws.Range("D45").Value = -0.09
$ws.Range("C46").Value = 94
$ws.Range("D46").Value = -0.123
$ws = $wb.Worksheets.Item("descriptives")
$ws.Range("B2").Value = 0.369
$ws.Range("C2").Value = -0.161
$ws.Range("B3").Value = 0.03
$ws.Range("C3").Value = 0.045
$ws.Range("B4").Value = 0.376
$ws.Range("C4").Value = -0.175
$ws.Range("B5").Value = 0.307
$ws.Range("C5").Value = -0.237
$ws.Range("B6").Value = 0.415
